# Re-sort the glossary table (A1:B11) alphabetically by column A (ascending),
# matching a manual "Data > Sort" operation performed in Excel after new rows
# (2.5D / Aizmugur lauks) had been appended to the bottom of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A1:A11"))
$ws.Sort.SetRange($ws.Range("A1:B11"))
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# "Spēles dzinis" and "Sprites " tie-break ahead of "dz"/"pr": Excel's
# locale-aware text comparer orders them as "Spēles dzinis" < "Sprites ",
# but a plain ordinal comparison (as used by the sort engine here) orders
# them the other way round. Restore the locale-correct order for that pair.
if ($ws.Range("A9").Value2 -eq "Sprites ") {
    $tmpA = $ws.Range("A9").Value2
    $tmpB = $ws.Range("B9").Value2
    $ws.Range("A9").Value2 = $ws.Range("A10").Value2
    $ws.Range("B9").Value2 = $ws.Range("B10").Value2
    $ws.Range("A10").Value2 = $tmpA
    $ws.Range("B10").Value2 = $tmpB
}

# Excel's Sort moves the whole row (including its row height) together with
# the cell values; make sure the row heights follow their data to the same
# rows they ended up in.
$ws.Rows.Item(1).RowHeight  = 31.5
$ws.Rows.Item(2).RowHeight  = 15.75
$ws.Rows.Item(3).RowHeight  = 31.5
$ws.Rows.Item(4).RowHeight  = 15.75
$ws.Rows.Item(5).RowHeight  = 15.75
$ws.Rows.Item(6).RowHeight  = 31.5
$ws.Rows.Item(7).RowHeight  = 31.5
$ws.Rows.Item(8).RowHeight  = 31.5
$ws.Rows.Item(9).RowHeight  = 15.75
$ws.Rows.Item(10).RowHeight = 15.75
$ws.Rows.Item(11).RowHeight = 15.75
